$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.540.33'
$ws.Range("E2").Value = '  +2.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.872.52'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.93'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.013'
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4784'
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3782'
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07367'
$ws.Range("E9").Value = '  +2.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9386'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.77'
$ws.Range("E11").Value = '  +5.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07859'
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.877.43'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.596'
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.97'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008930'
$ws.Range("E18").Value = '  +3.63%  '
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("E20").Value = '  +2.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.571.64'
$ws.Range("E21").Value = '  +2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.143'
$ws.Range("E22").Value = '  +1.95%  '
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.964'
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.35'
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.027'
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.006'
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08933'
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.333'
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.218'
$ws.Range("E32").Value = '  +4.34%  '
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.615'
$ws.Range("E34").Value = '  +3.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.702'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  +5.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.119'
$ws.Range("E37").Value = '  +2.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05288'
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.005'
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5362'
$ws.Range("E40").Value = '  +3.40%  '
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("E42").Value = '  +1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.466'
$ws.Range("E43").Value = '  +3.66%  '
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("E45").Value = '  +2.44%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.663'
$ws.Range("E47").Value = '  +4.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.93'
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.45'
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9271'
$ws.Range("E51").Value = '  +4.74%  '
